$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style captured from an untouched default-styled cell (B2),
# used to restore D-column cells to their original (unstyled) appearance
# after forcing a text NumberFormat so Excel does not auto-convert
# numeric-looking strings (e.g. "0.999") into real numbers.
$normalStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = '65.265.45'
$ws.Range("E2").Value = '  +4.64%  '
$ws.Range("D3").Value = '3.109.61'
$ws.Range("E3").Value = '  +2.54%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '561.54'
$ws.Range("D5").Style = $normalStyle
$ws.Range("E5").Value = '  +3.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.35'
$ws.Range("D6").Style = $normalStyle
$ws.Range("E6").Value = '  +7.31%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = $normalStyle
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '3.104.81'
$ws.Range("E8").Value = '  +2.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.501'
$ws.Range("D9").Style = $normalStyle
$ws.Range("E9").Value = '  +1.37%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.153'
$ws.Range("D10").Style = $normalStyle
$ws.Range("E10").Value = '  +3.18%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.31'
$ws.Range("D11").Style = $normalStyle
$ws.Range("E11").Value = '  +3.41%  '
$ws.Range("E12").Value = '  +4.37%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000231'
$ws.Range("D13").Style = $normalStyle
$ws.Range("E13").Value = '  +4.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.41'
$ws.Range("D14").Style = $normalStyle
$ws.Range("E14").Value = '  +2.28%  '
$ws.Range("D15").Value = '3.608.16'
$ws.Range("E15").Value = '  +2.23%  '
$ws.Range("D16").Value = '65.247.43'
$ws.Range("E16").Value = '  +4.55%  '
$ws.Range("D17").Value = '3.107.50'
$ws.Range("E17").Value = '  +2.32%  '
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.80'
$ws.Range("D19").Style = $normalStyle
$ws.Range("E19").Value = '  +1.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '480.73'
$ws.Range("D20").Style = $normalStyle
$ws.Range("E20").Value = '  -0.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.81'
$ws.Range("D21").Style = $normalStyle
$ws.Range("E21").Value = '  +4.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.680'
$ws.Range("D22").Style = $normalStyle
$ws.Range("E22").Value = '  +0.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.57'
$ws.Range("D23").Style = $normalStyle
$ws.Range("E23").Value = '  +6.98%  '
$ws.Range("E24").Value = '  +11.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.30'
$ws.Range("D25").Style = $normalStyle
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  +2.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.23'
$ws.Range("D28").Style = $normalStyle
$ws.Range("E28").Value = '  +5.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.07'
$ws.Range("D29").Style = $normalStyle
$ws.Range("E29").Value = '  +6.45%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.22'
$ws.Range("D31").Style = $normalStyle
$ws.Range("E31").Value = '  +1.59%  '
$ws.Range("E32").Value = '  +2.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.50'
$ws.Range("D33").Style = $normalStyle
$ws.Range("E33").Value = '  +5.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.65'
$ws.Range("D34").Style = $normalStyle
$ws.Range("E34").Value = '  -0.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.19'
$ws.Range("D35").Style = $normalStyle
$ws.Range("E35").Value = '  +5.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.10'
$ws.Range("D36").Style = $normalStyle
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '471.04'
$ws.Range("D37").Style = $normalStyle
$ws.Range("E37").Value = '  +3.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0412'
$ws.Range("D38").Style = $normalStyle
$ws.Range("E38").Value = '  +6.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0834'
$ws.Range("D39").Style = $normalStyle
$ws.Range("E39").Value = '  +3.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.92'
$ws.Range("D40").Style = $normalStyle
$ws.Range("E40").Value = '  +18.74%  '
$ws.Range("D41").Value = '3.008.64'
$ws.Range("E41").Value = '  -4.99%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.116'
$ws.Range("D42").Style = $normalStyle
$ws.Range("E42").Value = '  -1.57%  '
$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.27'
$ws.Range("D43").Style = $normalStyle
$ws.Range("E43").Value = '  +1.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '28.03'
$ws.Range("D44").Style = $normalStyle
$ws.Range("E44").Value = '  +6.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.261'
$ws.Range("D45").Style = $normalStyle
$ws.Range("E45").Value = '  +6.40%  '
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.15'
$ws.Range("D46").Style = $normalStyle
$ws.Range("E46").Value = '  +9.06%  '
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("D47").Style = $normalStyle
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.112'
$ws.Range("D48").Style = $normalStyle
$ws.Range("E48").Value = '  +3.04%  '
$ws.Range("D49").Value = '0.0₃0525'
$ws.Range("E49").Value = '  +5.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '115.17'
$ws.Range("D50").Style = $normalStyle
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("E51").Value = '  +2.12%  '
